$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell A1 from "ConversationID" to "id"
$ws.Range("A1").Value = "id"

# The header row was sized to wrap "ConversationID" onto two lines; with the
# shorter text "id" it no longer needs the extra height, so re-fit it back to
# the sheet's default single-line height.
$ws.Rows(1).AutoFit()

# Move the active selection to A2 (matches Excel's recorded selection after edit)
$ws.Range("A2").Select()
